function Set-CellText($Worksheet, $Cell, $Text) {
    $range = $Worksheet.Range($Cell)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "275.13"
Set-CellText $ws "G2" "4"
Set-CellText $ws "D3" "27.14"
Set-CellText $ws "E3" "1.92%"
Set-CellText $ws "G3" "4"
Set-CellText $ws "D4" "4.866"
Set-CellText $ws "E4" "-0.22%"
Set-CellText $ws "G4" "4"
Set-CellText $ws "D5" "0.06393"
Set-CellText $ws "E5" "1.26%"
Set-CellText $ws "G5" "4"
Set-CellText $ws "E6" "0.49%"
Set-CellText $ws "G6" "4"
Set-CellText $ws "D7" "1.209"
Set-CellText $ws "E7" "-2.89%"
Set-CellText $ws "G7" "4"
Set-CellText $ws "D8" "0.8769"
Set-CellText $ws "E8" "0.61%"
Set-CellText $ws "G8" "4"
Set-CellText $ws "D9" "0.1514"
Set-CellText $ws "E9" "3.99%"
Set-CellText $ws "G9" "4"
Set-CellText $ws "D10" "0.05113"
Set-CellText $ws "E10" "-0.65%"
Set-CellText $ws "G10" "4"
Set-CellText $ws "D11" "0.07522"
Set-CellText $ws "E11" "2.53%"
Set-CellText $ws "G11" "4"
Set-CellText $ws "D12" "0.02959"
Set-CellText $ws "E12" "-3.34%"
Set-CellText $ws "G12" "4"
Set-CellText $ws "D13" "0.08983"
Set-CellText $ws "E13" "-0.58%"
Set-CellText $ws "G13" "4"
Set-CellText $ws "D14" "0.001584"
Set-CellText $ws "E14" "0.80%"
Set-CellText $ws "G14" "4"
Set-CellText $ws "D15" "0.0006394"
Set-CellText $ws "E15" "1.24%"
Set-CellText $ws "G15" "4"
Set-CellText $ws "D16" "0.006185"
Set-CellText $ws "E16" "2.74%"
Set-CellText $ws "G16" "4"
Set-CellText $ws "D17" "3.468"
Set-CellText $ws "E17" "0.45%"
Set-CellText $ws "G17" "4"
Set-CellText $ws "D18" "3.314"
Set-CellText $ws "E18" "-1.15%"
Set-CellText $ws "G18" "4"
Set-CellText $ws "D19" "2.284"
Set-CellText $ws "E19" "0.01%"
Set-CellText $ws "G19" "4"
Set-CellText $ws "E20" "-0.94%"
Set-CellText $ws "G20" "4"
Set-CellText $ws "D21" "0.1349"
Set-CellText $ws "E21" "1.88%"
Set-CellText $ws "G21" "4"
Set-CellText $ws "D22" "3.904"
Set-CellText $ws "E22" "-0.03%"
Set-CellText $ws "G22" "4"
Set-CellText $ws "D23" "0.04428"
Set-CellText $ws "E23" "0.16%"
Set-CellText $ws "G23" "4"
Set-CellText $ws "G24" "4"
Set-CellText $ws "E25" "-0.10%"
Set-CellText $ws "G25" "4"
Set-CellText $ws "D26" "0.003853"
Set-CellText $ws "E26" "-12.57%"
Set-CellText $ws "G26" "4"
Set-CellText $ws "D27" "0.0001200"
Set-CellText $ws "E27" "-0.04%"
Set-CellText $ws "G27" "4"
Set-CellText $ws "E28" "14.16%"
Set-CellText $ws "G28" "4"
Set-CellText $ws "G29" "4"
Set-CellText $ws "G30" "4"
Set-CellText $ws "G31" "4"
Set-CellText $ws "G32" "4"
Set-CellText $ws "G33" "4"
Set-CellText $ws "G34" "4"
Set-CellText $ws "G35" "4"
Set-CellText $ws "G36" "4"
Set-CellText $ws "G37" "4"
Set-CellText $ws "G38" "4"
Set-CellText $ws "G39" "4"
Set-CellText $ws "E40" "2.58%"
Set-CellText $ws "G40" "4"
Set-CellText $ws "D41" "0.006808"
Set-CellText $ws "E41" "1.04%"
Set-CellText $ws "G41" "4"
Set-CellText $ws "E42" "0.74%"
Set-CellText $ws "G42" "4"
Set-CellText $ws "D43" "0.002101"
Set-CellText $ws "E43" "-0.52%"
Set-CellText $ws "G43" "4"
Set-CellText $ws "D44" "0.01149"
Set-CellText $ws "E44" "-8.35%"
Set-CellText $ws "G44" "4"
Set-CellText $ws "D45" "0.00005213"
Set-CellText $ws "E45" "-1.94%"
Set-CellText $ws "G45" "4"
Set-CellText $ws "E46" "-45.57%"
Set-CellText $ws "G46" "4"
Set-CellText $ws "E47" "0.13%"
Set-CellText $ws "G47" "4"
Set-CellText $ws "G48" "4"
Set-CellText $ws "G49" "4"
Set-CellText $ws "G50" "4"
Set-CellText $ws "G51" "4"
